$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header cell B1
$ws.Range("B1").Value = "ishan"

# Add new header cell C1, copying the style used by the other header cells (A1/B1)
$ws.Range("C1").Value = "Vikas"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in column C data for existing rows
$ws.Range("C2").Value = "Eligible"
$ws.Range("C3").Value = "Not Eligible"

# Add new row 4
$ws.Range("A4").Value = "Is elegible"
$ws.Range("B4").Value = "Eligible"
$ws.Range("C4").Value = "Not Eligible"

# Copy formatting from A3 (header-style column A cell) to A4
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
